$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column G: "Number of lines" per contract file
$ws.Range("G1").Value = "Number of lines"

$lineCounts = @(
    284, 703, 134, 171, 46, 103, 84, 74, 410, 50, 6, 23,
    68, 15, 26, 3, 23, 3, 88, 3, 21, 5, 113, 249,
    98, 29, 225, 457, 471, 50, 98, 364, 11, 197, 97, 344,
    235, 519, 66, 317, 86, 91, 7, 139, 155, 43, 203, 14,
    103, 5, 3, 23, 3, 76, 3, 40, 33, 5, 147, 283,
    87, 28, 170, 443, 171, 461, 45, 194, 86, 376, 237, 477,
    56, 212, 1019, 247
)

for ($i = 0; $i -lt $lineCounts.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $lineCounts[$i]
}

